# Apply "commit before going home" edit:
#  - Add a "routers" / "ejs page" header pair in F3/G3
#  - Mark "v" in F4 and F9 (rows whose router + ejs page have been checked)
#  - Selection moves to D20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F4").Value = "v"
$ws.Range("F9").Value = "v"

$ws.Range("F3").Value = "routers"
$ws.Range("G3").Value = "ejs page"

$ws.Range("D20").Select()
